$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.057.99"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "2.512.46"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "532.03"
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.11"
$ws.Range("E6").Value = "  -3.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.564"
$ws.Range("E8").Value = "  -1.01%  "
$ws.Range("D9").Value = "2.517.55"
$ws.Range("E9").Value = "  -0.44%  "
$ws.Range("E10").Value = "  -0.14%  "
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.47"
$ws.Range("E12").Value = "  -2.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.356"
$ws.Range("E13").Value = "  +0.98%  "
$ws.Range("D14").Value = "2.959.28"
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.07"
$ws.Range("E15").Value = "  -1.84%  "
$ws.Range("D16").Value = "59.022.26"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").Value = "2.513.55"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.99"
$ws.Range("E19").Value = "  -1.92%  "
$ws.Range("E20").Value = "  -0.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.85"
$ws.Range("E21").Value = "  -0.88%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.83"
$ws.Range("E23").Value = "  +0.79%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.25"
$ws.Range("E24").Value = "  +0.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.423"
$ws.Range("E25").Value = "  -3.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.165"
$ws.Range("E26").Value = "  +1.87%  "
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.80"
$ws.Range("E28").Value = "  +0.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.74"
$ws.Range("E29").Value = "  +0.99%  "
$ws.Range("D30").Value = "0.0₃0770"
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("E31").Value = "  -0.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "161.87"
$ws.Range("E32").Value = "  +3.22%  "
$ws.Range("E33").Value = "  +0.26%  "
$ws.Range("E34").Value = "  -6.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.44"
$ws.Range("E35").Value = "  +0.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.45"
$ws.Range("E36").Value = "  -1.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.23"
$ws.Range("E37").Value = "  -2.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.58"
$ws.Range("E38").Value = "  -1.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.96"
$ws.Range("E39").Value = "  +0.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.64"
$ws.Range("E40").Value = "  -1.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.804"
$ws.Range("E41").Value = "  -2.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.23"
$ws.Range("E42").Value = "  -7.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "279.71"
$ws.Range("E43").Value = "  -6.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  +0.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.85"
$ws.Range("E45").Value = "  +0.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.593"
$ws.Range("E46").Value = "  -1.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0929"
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.09"
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.34"
$ws.Range("E49").Value = "  -1.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0509"
$ws.Range("E50").Value = "  -0.97%  "
$ws.Range("E51").Value = "  -2.34%  "
